$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold, borders, centered) from H1 to I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows: I and J columns (rows 2-13)
$values = @(
    @(8, 9),   # row 2
    @(8, 9),   # row 3
    @(8, 8),   # row 4
    @(4, 5),   # row 5
    @(4, 6),   # row 6
    @(6, 6),   # row 7
    @(6, 7),   # row 8
    @(8, 8),   # row 9
    @(4, 4),   # row 10
    @(6, 6),   # row 11
    @(7, 7),   # row 12
    @(8, 8)    # row 13
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
